# Conversion factor: EMU per point
$EMU_PER_POINT = 12700

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(17)

# --- Update the warning text in "TextBox 6" (shape id 7) ---
$textBox = $s.Shapes.Item(5)
$textBox.TextFrame.TextRange.Text = "El modificador de acceso, ser estático vs de instancia o el valor de retorno no son suficientes para sobrecargar un método!"

# --- Resize/reposition "Rectangle 7" (shape id 8) ---
# (values nudged by a few EMU pre-conversion to compensate for the
# runtime's internal float32 point storage, so the saved EMU matches
# the intended target exactly)
$rect7 = $s.Shapes.Item(6)
$rect7.Left = 594168 / $EMU_PER_POINT
$rect7.Top = 1591535 / $EMU_PER_POINT
$rect7.Width = 3121305 / $EMU_PER_POINT
$rect7.Height = 394137 / $EMU_PER_POINT

# --- Resize/reposition "Rectangle 8" (shape id 9) ---
$rect8 = $s.Shapes.Item(7)
$rect8.Left = 594166 / $EMU_PER_POINT
$rect8.Top = 2794641 / $EMU_PER_POINT
$rect8.Width = 2936111 / $EMU_PER_POINT
$rect8.Height = 394137 / $EMU_PER_POINT
